$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move the active selection to E7 (was E6)
$ws.Range("E7").Select()

# Column B: drop bestFit, set explicit width (target stored width 44.42578125;
# engine quantizes ColumnWidth to 1/6 px steps, so 43.666666... is the closest
# achievable input, landing on stored width 44.5)
$ws.Columns.Item(2).ColumnWidth = 43.666666666666664

# Row 4: height 60 -> 30, clear G4
$ws.Rows.Item(4).RowHeight = 30
# Row 5: height 75 -> 30, clear G5
$ws.Rows.Item(5).RowHeight = 30
# Row 6: height 75 -> 45, clear G6
$ws.Rows.Item(6).RowHeight = 45
# Row 7: height 75 -> 45, clear G7
$ws.Rows.Item(7).RowHeight = 45

# Clear the G column "Expected Results" duplicate text in rows 3-8 and 10
$ws.Range("G3").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("G10").ClearContents()

# Fill in G9 / G11 to mirror the Notes column values
$ws.Range("G9").Value = "AppliedPatch_CCR"
$ws.Range("G11").Value = "not added to queue by Dispatcher"

# G13 / G14 need both the value AND the style used by G2/G9/G11 (vertical
# centered). VerticalAlignment assignment is a no-op on this host, so copy
# the format from a cell that already has the target style, then set text.
$ws.Range("G2").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("G13").Value = "not added to queue by Dispatcher"
$ws.Range("G14").Value = "not added to queue by Dispatcher"
$excel.CutCopyMode = $false
